$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 rows at the top
$ws.Rows("1:5").Insert()

$ws.Range("B1").Value = "OPERADORES S2S - CSV MERA"
$ws.Range("B1:G4").Merge()
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Color = 16777215
$ws.Range("B1").Font.Size = 22
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1:G4").Interior.Color = 6968388
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("A6:H6").Interior.Color = 6968388
$ws.Range("A6:H6").Font.Bold = $true
$ws.Range("A6:H6").Font.Color = 16777215
